# Re-generate the quadratic/linear problem data (new random draw of the
# generator), per commit "volver a generar problemas cuadraticos y lineales".
#
# Only the *values* shown by the generator change; headers, expression
# labels that stay the same (e.g. Restriction_Set_Type) and sheet layout
# are untouched.
#
# Because several of these columns hold numbers formatted as genuine TEXT
# (shared-string) cells -- e.g. "-1.1500000000000004" -- rather than native
# numeric cells, we force the destination range to Text format before
# writing so Excel doesn't silently re-interpret the literal as a number,
# then restore the cell style back to Normal so we don't leave behind a
# stray number-format override.
function Set-TextValues {
    param(
        $Range,
        [object[]]$Values
    )
    $Range.NumberFormat = "@"
    $Range.Value = $Values
    $Range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Restricciones_del_follower: expression shifted + new Lambda/Beta/Gamma
# ---------------------------------------------------------------------
$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")

Set-TextValues $wsFollower.Range("A2") "-12.85 + x + 2y"
Set-TextValues $wsFollower.Range("B2") "-1.1500000000000004"
Set-TextValues $wsFollower.Range("D2") "0.07"
Set-TextValues $wsFollower.Range("E2") "0"
Set-TextValues $wsFollower.Range("F2") "1.0"

Set-TextValues $wsFollower.Range("A3") "3.3499999999999996 + x - 2y"
Set-TextValues $wsFollower.Range("B3") "-5.35"
Set-TextValues $wsFollower.Range("D3") "0.21"
Set-TextValues $wsFollower.Range("E3") "6.0"
Set-TextValues $wsFollower.Range("F3") "5.8"

Set-TextValues $wsFollower.Range("A4") "-7.45 - 2x + y"
Set-TextValues $wsFollower.Range("B4") "-6.45"
Set-TextValues $wsFollower.Range("D4") "0.4"
Set-TextValues $wsFollower.Range("E4") "1.7000000000000002"
Set-TextValues $wsFollower.Range("F4") "1.3"

# ---------------------------------------------------------------------
# Punto_modificado: new (x, y) point
# ---------------------------------------------------------------------
$wsPunto = $wb.Worksheets.Item("Punto_modificado")
Set-TextValues $wsPunto.Range("A2") "4.75"
Set-TextValues $wsPunto.Range("B2") "4.05"

# ---------------------------------------------------------------------
# Vector_bf: new value
# ---------------------------------------------------------------------
$wsBf = $wb.Worksheets.Item("Vector_bf")
Set-TextValues $wsBf.Range("A2") "1.7800000000000002"

# ---------------------------------------------------------------------
# Vector_BF: new values
# ---------------------------------------------------------------------
$wsBF = $wb.Worksheets.Item("Vector_BF")
Set-TextValues $wsBF.Range("A2") "-6.1"
Set-TextValues $wsBF.Range("A3") "6.200000000000001"

Write-Output "done"
